$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header / value for column G ("adsriver" source tracking column)
$ws.Range("G1").Value = "adsriver"

# Mark "ok" for most rows (G3 and G7 are intentionally left blank)
$ws.Range("G2").Value = "ok"
$ws.Range("G4").Value = "ok"
$ws.Range("G5").Value = "ok"
$ws.Range("G6").Value = "ok"
$ws.Range("G8").Value = "ok"
$ws.Range("G9").Value = "ok"
$ws.Range("G10").Value = "ok"
$ws.Range("G11").Value = "ok"
$ws.Range("G12").Value = "ok"
$ws.Range("G13").Value = "ok"
$ws.Range("G14").Value = "ok"
$ws.Range("G15").Value = "ok"
$ws.Range("G16").Value = "ok"
$ws.Range("G17").Value = "ok"

# Final selection matches the author's last edited cell
$ws.Range("G4").Select()
